$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Photo Links")

$rows = @(8, 9, 10, 12, 13, 14, 17, 18, 20, 21, 23, 24, 26, 27, 28, 39, 42)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -like "*.png") {
        $newVal = $val.Substring(0, $val.Length - 4) + ".jpg"
        $cell.Value = $newVal
    }
}
